$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column T ("metaParamTex"): the LaTeX symbol of each distribution's
# meta-parameter. This mirrors the existing "paramTex" column (S) but groups
# the Bernoulli-family rows under \pi and keeps \beta / \lambda for the rest.
$ws.Range("T1").Value  = "metaParamTex"
$ws.Range("T2").Value  = "\pi"
$ws.Range("T3").Value  = "\pi"
$ws.Range("T4").Value  = "\pi"
$ws.Range("T5").Value  = "\beta"
$ws.Range("T6").Value  = "\beta"
$ws.Range("T7").Value  = "\beta"
$ws.Range("T8").Value  = "\beta"
$ws.Range("T9").Value  = "\lambda"
$ws.Range("T10").Value = "\lambda"
$ws.Range("T11").Value = "\lambda"
$ws.Range("T12").Value = "\lambda"
$ws.Range("T13").Value = "\lambda"
$ws.Range("T14").Value = "\lambda"

# Move the live selection onto the newly-filled range.
$ws.Range("T9:T14").Select()
